$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 20142.715
$ws.Range("I2").Value = 199.8
$ws.Range("K2").Value = 199.8
$ws.Range("M2").Value = -86.80000000000001

$ws.Range("H5").Value = 197.28572
$ws.Range("I5").Value = 197.28572
$ws.Range("K5").Value = 197.28572
$ws.Range("M5").Value = -82.28572

$ws.Range("H9").Value = 11216.667
$ws.Range("I9").Value = 12618.5
$ws.Range("K9").Value = 12618.5
$ws.Range("M9").Value = -12449.5

$ws.Range("H12").Value = 16894.334
$ws.Range("I12").Value = 20202
$ws.Range("K12").Value = 20202
$ws.Range("M12").Value = -20032

$ws.Range("H17").Value = 643271.2
$ws.Range("J17").Value = 763832
$ws.Range("L17").Value = 2291496
$ws.Range("N17").Value = -2291832

$ws.Range("H19").Value = 1078.3
$ws.Range("I19").Value = 1119.8
$ws.Range("K19").Value = 1119.8
$ws.Range("M19").Value = -944.8

$ws.Range("H33").Value = 268.3846
$ws.Range("J33").Value = 320
$ws.Range("L33").Value = 320
$ws.Range("N33").Value = -778

$ws.Range("H39").Value = 83334856
$ws.Range("J39").Value = 1923.2
$ws.Range("L39").Value = 5769.6
$ws.Range("N39").Value = -6361.6

$ws.Range("H41").Value = 1482.6666
$ws.Range("I41").Value = 1763.8572
$ws.Range("K41").Value = 1763.8572
$ws.Range("M41").Value = -1323.8572

$ws.Range("H42").Value = 166670480
$ws.Range("I42").Value = 250000220
$ws.Range("K42").Value = 750000660
$ws.Range("M42").Value = -750000430

$ws.Range("H43").Value = 13571.429
$ws.Range("I43").Value = 8999.75
$ws.Range("J43").Value = 19667
$ws.Range("K43").Value = 8999.75
$ws.Range("L43").Value = 19667
$ws.Range("M43").Value = -8930.75
$ws.Range("N43").Value = -19805

$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

$ws.Range("H60").Value = 4000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null

$ws.Range("H64").Value = 4328.6113
$ws.Range("I64").Value = 4079.5557
$ws.Range("J64").Value = 4577.6665
$ws.Range("K64").Value = 4079.5557
$ws.Range("L64").Value = 4577.6665
$ws.Range("M64").Value = -3831.5557
$ws.Range("N64").Value = -5073.6665

$ws.Range("H67").Value = 4328.6113
$ws.Range("I67").Value = 4079.5557
$ws.Range("J67").Value = 4577.6665
$ws.Range("K67").Value = 4079.5557
$ws.Range("L67").Value = 4577.6665
$ws.Range("M67").Value = -3221.5557
$ws.Range("N67").Value = -6293.6665

$ws.Range("H70").Value = 5144.125
$ws.Range("I70").Value = 1199.6666
$ws.Range("J70").Value = 7510.8
$ws.Range("K70").Value = 3598.9998
$ws.Range("L70").Value = 22532.4
$ws.Range("M70").Value = -3328.9998
$ws.Range("N70").Value = -23072.4

$ws.Range("H73").Value = 5144.125
$ws.Range("I73").Value = 1199.6666
$ws.Range("J73").Value = 7510.8
$ws.Range("K73").Value = 3598.9998
$ws.Range("L73").Value = 22532.4
$ws.Range("M73").Value = -2662.9998
$ws.Range("N73").Value = -24404.4

$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872

$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360

$ws.Range("H94").Value = 2506.8333
$ws.Range("I94").Value = 2468.4
$ws.Range("K94").Value = 2468.4
$ws.Range("M94").Value = -2017.4

$ws.Range("H98").Value = 1328.5
$ws.Range("I98").Value = 1485.7142
$ws.Range("K98").Value = 1485.7142
$ws.Range("M98").Value = 12.28580000000011

$ws.Range("H107").Value = 1986.4286
$ws.Range("I107").Value = 1408.4615
$ws.Range("K107").Value = 1408.4615
$ws.Range("M107").Value = 511.5385000000001

$ws.Range("H112").Value = 334349
$ws.Range("I112").Value = 201219
$ws.Range("J112").Value = 999999
$ws.Range("K112").Value = 603657
$ws.Range("L112").Value = 2999997
$ws.Range("M112").Value = -602549
$ws.Range("N112").Value = -3002213

$ws.Range("H113").Value = 2535
$ws.Range("I113").Value = 2611.5715
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 2611.5715
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 642.4285
$ws.Range("N113").Value = -8507

$ws.Range("H114").Value = 69999
$ws.Range("J114").Value = 69999
$ws.Range("L114").Value = 69999
$ws.Range("N114").Value = -78677

$ws.Range("H116").Value = 4571.3687
$ws.Range("I116").Value = 4361.5713
$ws.Range("K116").Value = 4361.5713
$ws.Range("M116").Value = -919.5712999999996

$ws.Range("H122").Value = 1328.5
$ws.Range("I122").Value = 1485.7142
$ws.Range("K122").Value = 4457.142599999999
$ws.Range("M122").Value = -2007.142599999999

$ws.Range("H125").Value = 17745656
$ws.Range("I125").Value = 4239463
$ws.Range("J125").Value = 31251850
$ws.Range("K125").Value = 38155167
$ws.Range("L125").Value = 281266650
$ws.Range("M125").Value = -38152707
$ws.Range("N125").Value = -281271570

$ws.Range("H132").Value = 3200.3333
$ws.Range("I132").Value = 3285.1924
$ws.Range("K132").Value = 9855.5772
$ws.Range("M132").Value = -7325.5772

$ws.Range("H135").Value = 88235750
$ws.Range("I135").Value = 38461944
$ws.Range("J135").Value = 250000660
$ws.Range("K135").Value = 346157496
$ws.Range("L135").Value = 2250005940
$ws.Range("M135").Value = -346154961
$ws.Range("N135").Value = -2250011010

$ws.Range("H137").Value = 2227.5386
$ws.Range("I137").Value = 1541
$ws.Range("J137").Value = 3028.5
$ws.Range("K137").Value = 4623
$ws.Range("L137").Value = 9085.5
$ws.Range("M137").Value = -2073
$ws.Range("N137").Value = -14185.5

$ws.Range("H138").Value = 3950.652
$ws.Range("I138").Value = 2675.0908
$ws.Range("J138").Value = 4351.543
$ws.Range("K138").Value = 8025.2724
$ws.Range("L138").Value = 13054.629
$ws.Range("M138").Value = -2885.2724
$ws.Range("N138").Value = -23334.629

$ws.Range("H141").Value = 4639.5713
$ws.Range("I141").Value = 3747
$ws.Range("K141").Value = 11241
$ws.Range("M141").Value = -6061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2452103.2
$ws.Range("I2").Value = 2941724
$ws.Range("K2").Value = 2941724
$ws.Range("M2").Value = -2941611

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = $null
$ws.Range("N7").Value = $null

$ws.Range("H32").Value = 3229.7458
$ws.Range("I32").Value = 1704.7407
$ws.Range("K32").Value = 1704.7407
$ws.Range("M32").Value = -1417.7407

$ws.Range("H33").Value = 4000
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3671

$ws.Range("H39").Value = 1520.5
$ws.Range("I39").Value = 1520.5
$ws.Range("K39").Value = 1520.5
$ws.Range("M39").Value = -1000.5

$ws.Range("H45").Value = 3783.7144
$ws.Range("I45").Value = 4314.375
$ws.Range("K45").Value = 4314.375
$ws.Range("M45").Value = -3937.375

$ws.Range("H61").Value = 50002700
$ws.Range("I61").Value = 58825496
$ws.Range("K61").Value = 58825496
$ws.Range("M61").Value = -58825284

$ws.Range("H74").Value = 32259726
$ws.Range("I74").Value = 37038504
$ws.Range("J74").Value = 2975.25
$ws.Range("K74").Value = 37038504
$ws.Range("L74").Value = 2975.25
$ws.Range("M74").Value = -37037630
$ws.Range("N74").Value = -4723.25

$ws.Range("H77").Value = 32259726
$ws.Range("I77").Value = 37038504
$ws.Range("J77").Value = 2975.25
$ws.Range("K77").Value = 185192520
$ws.Range("L77").Value = 14876.25
$ws.Range("M77").Value = -185188152
$ws.Range("N77").Value = -23612.25

$ws.Range("H88").Value = 2330.6667
$ws.Range("I88").Value = 2119
$ws.Range("K88").Value = 2119
$ws.Range("M88").Value = -1713

$ws.Range("H91").Value = 2330.6667
$ws.Range("I91").Value = 2119
$ws.Range("K91").Value = 2119
$ws.Range("M91").Value = -715

$ws.Range("H96").Value = 24995
$ws.Range("J96").Value = 24995
$ws.Range("L96").Value = 24995
$ws.Range("N96").Value = -30487

$ws.Range("H97").Value = 497.85715
$ws.Range("I97").Value = 506.66666
$ws.Range("J97").Value = 445
$ws.Range("K97").Value = 506.66666
$ws.Range("L97").Value = 445
$ws.Range("M97").Value = -10.66665999999998
$ws.Range("N97").Value = -1437

$ws.Range("H102").Value = 9096876
$ws.Range("I102").Value = 10006164
$ws.Range("K102").Value = 10006164
$ws.Range("M102").Value = -10004542

$ws.Range("H110").Value = 59931.234
$ws.Range("I110").Value = 63614.438
$ws.Range("K110").Value = 63614.438
$ws.Range("M110").Value = -61569.438

$ws.Range("H116").Value = 2452103.2
$ws.Range("I116").Value = 2941724
$ws.Range("K116").Value = 2941724
$ws.Range("M116").Value = -2939430

$ws.Range("H122").Value = 4902.6523
$ws.Range("I122").Value = 1888.3
$ws.Range("K122").Value = 5664.9
$ws.Range("M122").Value = -3214.9

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = $null
$ws.Range("N129").Value = $null

$ws.Range("H132").Value = 5003239
$ws.Range("I132").Value = 7146149
$ws.Range("K132").Value = 21438447
$ws.Range("M132").Value = -21435917

$ws.Range("H136").Value = 50002700
$ws.Range("I136").Value = 58825496
$ws.Range("K136").Value = 176476488
$ws.Range("M136").Value = -176473938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2452103.2
$ws.Range("I3").Value = 2941724
$ws.Range("K3").Value = 2941724
$ws.Range("M3").Value = -2941610

$ws.Range("H20").Value = 2470.1
$ws.Range("I20").Value = 2344.0715
$ws.Range("J20").Value = 2764.1667
$ws.Range("K20").Value = 2344.0715
$ws.Range("L20").Value = 2764.1667
$ws.Range("M20").Value = -2097.0715
$ws.Range("N20").Value = -3258.1667

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = $null
$ws.Range("N75").Value = $null

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = $null
$ws.Range("N78").Value = $null

$ws.Range("H94").Value = 1282.3529
$ws.Range("I94").Value = 1307.0667
$ws.Range("J94").Value = 1097
$ws.Range("K94").Value = 1307.0667
$ws.Range("L94").Value = 1097
$ws.Range("M94").Value = -856.0667000000001
$ws.Range("N94").Value = -1999

$ws.Range("H105").Value = 3292.4375
$ws.Range("I105").Value = 3098.1667
$ws.Range("K105").Value = 3098.1667
$ws.Range("M105").Value = -1351.1667

$ws.Range("H134").Value = 28337872
$ws.Range("I134").Value = 28337872
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 85013616
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -85011081
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4003200
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 20000000
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 20000000
$ws.Range("M4").Value = -3888
$ws.Range("N4").Value = -20000224

$ws.Range("H7").Value = 6600.5557
$ws.Range("I7").Value = 9077.691999999999
$ws.Range("K7").Value = 9077.691999999999
$ws.Range("M7").Value = -8964.691999999999

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = $null
$ws.Range("N17").Value = $null

$ws.Range("H22").Value = 489.47223
$ws.Range("I22").Value = 497.74286
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 497.74286
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -147.74286
$ws.Range("N22").Value = -900

$ws.Range("H25").Value = 7000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 7000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 7000
$ws.Range("N25").Value = -7348
$ws.Range("M25").Value = $null

$ws.Range("H31").Value = 3608.3215
$ws.Range("I31").Value = 3840.0435
$ws.Range("J31").Value = 2542.4
$ws.Range("K31").Value = 3840.0435
$ws.Range("L31").Value = 2542.4
$ws.Range("M31").Value = -3545.0435
$ws.Range("N31").Value = -3132.4

$ws.Range("H34").Value = 3608.3215
$ws.Range("I34").Value = 3840.0435
$ws.Range("J34").Value = 2542.4
$ws.Range("K34").Value = 3840.0435
$ws.Range("L34").Value = 2542.4
$ws.Range("M34").Value = -3638.0435
$ws.Range("N34").Value = -2946.4

$ws.Range("H58").Value = 13896047
$ws.Range("I58").Value = 17249046
$ws.Range("J58").Value = 5052.4287
$ws.Range("K58").Value = 17249046
$ws.Range("L58").Value = 5052.4287
$ws.Range("M58").Value = -17248843
$ws.Range("N58").Value = -5458.4287

$ws.Range("H74").Value = 39593.6
$ws.Range("J74").Value = 39593.6
$ws.Range("L74").Value = 39593.6
$ws.Range("N74").Value = -41341.6

$ws.Range("H77").Value = 39593.6
$ws.Range("J77").Value = 39593.6
$ws.Range("L77").Value = 118780.8
$ws.Range("N77").Value = -127516.8

$ws.Range("H134").Value = 25002336
$ws.Range("I134").Value = 25002336
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 75007008
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -75004473
$ws.Range("N134").Value = $null

$ws.Range("H136").Value = 13896047
$ws.Range("I136").Value = 17249046
$ws.Range("J136").Value = 5052.4287
$ws.Range("K136").Value = 51747138
$ws.Range("L136").Value = 15157.2861
$ws.Range("M136").Value = -51744588
$ws.Range("N136").Value = -20257.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 407207.72
$ws.Range("I4").Value = 482079.97
$ws.Range("J4").Value = 2897.6
$ws.Range("K4").Value = 1446239.91
$ws.Range("L4").Value = 8692.799999999999
$ws.Range("M4").Value = -1446127.91
$ws.Range("N4").Value = -8916.799999999999

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null

$ws.Range("H11").Value = 120341.93
$ws.Range("I11").Value = 130933
$ws.Range("J11").Value = 51500
$ws.Range("K11").Value = 392799
$ws.Range("L11").Value = 154500
$ws.Range("M11").Value = -392659
$ws.Range("N11").Value = -154780

$ws.Range("H37").Value = 114464.07
$ws.Range("J37").Value = 114464.07
$ws.Range("L37").Value = 343392.21
$ws.Range("N37").Value = -343616.21

$ws.Range("H38").Value = 181.3158
$ws.Range("J38").Value = 173.57143
$ws.Range("L38").Value = 520.71429
$ws.Range("N38").Value = -1214.71429

$ws.Range("H47").Value = 1126.8334
$ws.Range("I47").Value = 1126.8334
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 3380.5002
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -2949.5002
$ws.Range("N47").Value = $null

$ws.Range("H50").Value = 1350.7778
$ws.Range("I50").Value = 878.7143
$ws.Range("K50").Value = 2636.1429
$ws.Range("M50").Value = -2155.1429

$ws.Range("H51").Value = 1681.6666
$ws.Range("I51").Value = 1681.6666
$ws.Range("K51").Value = 5044.9998
$ws.Range("M51").Value = -4584.9998

$ws.Range("H53").Value = 1350.7778
$ws.Range("I53").Value = 878.7143
$ws.Range("K53").Value = 2636.1429
$ws.Range("M53").Value = -2155.1429

$ws.Range("H68").Value = 2902.2856
$ws.Range("J68").Value = 2992.8462
$ws.Range("L68").Value = 8978.5386
$ws.Range("N68").Value = -10600.5386

$ws.Range("H70").Value = 15971
$ws.Range("I70").Value = 9951.666999999999
$ws.Range("J70").Value = 25000
$ws.Range("K70").Value = 29855.001
$ws.Range("L70").Value = 75000
$ws.Range("M70").Value = -29540.001
$ws.Range("N70").Value = -75630

$ws.Range("H71").Value = 2902.2856
$ws.Range("J71").Value = 2992.8462
$ws.Range("L71").Value = 26935.6158
$ws.Range("N71").Value = -35047.6158

$ws.Range("H73").Value = 15971
$ws.Range("I73").Value = 9951.666999999999
$ws.Range("J73").Value = 25000
$ws.Range("K73").Value = 29855.001
$ws.Range("L73").Value = 75000
$ws.Range("M73").Value = -28763.001
$ws.Range("N73").Value = -77184

$ws.Range("H75").Value = 4938.3335
$ws.Range("J75").Value = 4907.5
$ws.Range("L75").Value = 14722.5
$ws.Range("N75").Value = -16718.5

$ws.Range("H78").Value = 4938.3335
$ws.Range("J78").Value = 4907.5
$ws.Range("L78").Value = 44167.5
$ws.Range("N78").Value = -54151.5

$ws.Range("H113").Value = 67243.336
$ws.Range("I113").Value = 125288.25
$ws.Range("K113").Value = 375864.75
$ws.Range("M113").Value = -373694.75

$ws.Range("H137").Value = 2992.5
$ws.Range("J137").Value = 3500
$ws.Range("L137").Value = 10500
$ws.Range("N137").Value = -20700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 45172.5
$ws.Range("J34").Value = 45172.5
$ws.Range("L34").Value = 45172.5
$ws.Range("N34").Value = -45708.5

$ws.Range("H64").Value = 73500
$ws.Range("J64").Value = 73500
$ws.Range("L64").Value = 73500
$ws.Range("N64").Value = -73996

$ws.Range("H67").Value = 73500
$ws.Range("J67").Value = 73500
$ws.Range("L67").Value = 73500
$ws.Range("N67").Value = -75216

$ws.Range("H70").Value = 6182.5
$ws.Range("I70").Value = 5524.25
$ws.Range("K70").Value = 5524.25
$ws.Range("M70").Value = -5254.25

$ws.Range("H73").Value = 6182.5
$ws.Range("I73").Value = 5524.25
$ws.Range("K73").Value = 5524.25
$ws.Range("M73").Value = -4588.25

$ws.Range("H76").Value = 45172.5
$ws.Range("J76").Value = 45172.5
$ws.Range("L76").Value = 45172.5
$ws.Range("N76").Value = -45802.5

$ws.Range("H79").Value = 45172.5
$ws.Range("J79").Value = 45172.5
$ws.Range("L79").Value = 45172.5
$ws.Range("N79").Value = -47356.5

$ws.Range("H94").Value = 7000
$ws.Range("I94").Value = 5000
$ws.Range("K94").Value = 5000
$ws.Range("M94").Value = -4324

$ws.Range("H99").Value = 24799.334
$ws.Range("I99").Value = 2199.5
$ws.Range("K99").Value = 2199.5
$ws.Range("M99").Value = 46.5

$ws.Range("H102").Value = 1903.72
$ws.Range("I102").Value = 1920.8334
$ws.Range("K102").Value = 1920.8334
$ws.Range("M102").Value = -298.8334

$ws.Range("H104").Value = 64402.2
$ws.Range("J104").Value = 64402.2
$ws.Range("L104").Value = 64402.2
$ws.Range("N104").Value = -71390.2

$ws.Range("H113").Value = 103419.3
$ws.Range("I113").Value = 146528.28
$ws.Range("J113").Value = 2831.6667
$ws.Range("K113").Value = 146528.28
$ws.Range("L113").Value = 2831.6667
$ws.Range("M113").Value = -144358.28
$ws.Range("N113").Value = -7171.6667

$ws.Range("H122").Value = 3112.647
$ws.Range("I122").Value = 2075.2
$ws.Range("J122").Value = 5994.4443
$ws.Range("K122").Value = 6225.599999999999
$ws.Range("L122").Value = 17983.3329
$ws.Range("M122").Value = -3775.599999999999
$ws.Range("N122").Value = -22883.3329

$ws.Range("H132").Value = 17859676
$ws.Range("I132").Value = 17859676
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 53579028
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -53576498
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 77798.8
$ws.Range("J136").Value = 77798.8
$ws.Range("L136").Value = 233396.4
$ws.Range("N136").Value = -238496.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5062.6665
$ws.Range("I7").Value = 4875.2
$ws.Range("K7").Value = 4875.2
$ws.Range("M7").Value = -4763.2

$ws.Range("H22").Value = 1859.3334
$ws.Range("I22").Value = 2138.3076
$ws.Range("J22").Value = 1406
$ws.Range("K22").Value = 2138.3076
$ws.Range("L22").Value = 1406
$ws.Range("M22").Value = -1843.3076
$ws.Range("N22").Value = -1996

$ws.Range("H27").Value = 1859.3334
$ws.Range("I27").Value = 2138.3076
$ws.Range("J27").Value = 1406
$ws.Range("K27").Value = 2138.3076
$ws.Range("L27").Value = 1406
$ws.Range("M27").Value = -2031.3076
$ws.Range("N27").Value = -1620

$ws.Range("H40").Value = 3018.375
$ws.Range("I40").Value = 3018.375
$ws.Range("K40").Value = 3018.375
$ws.Range("M40").Value = -2882.375

$ws.Range("H46").Value = 1991.4706
$ws.Range("I46").Value = 2061
$ws.Range("K46").Value = 2061
$ws.Range("M46").Value = -1873

$ws.Range("H68").Value = 3716.6365
$ws.Range("I68").Value = 1870.125
$ws.Range("J68").Value = 8640.666999999999
$ws.Range("K68").Value = 1870.125
$ws.Range("L68").Value = 8640.666999999999
$ws.Range("M68").Value = -1121.125
$ws.Range("N68").Value = -10138.667

$ws.Range("H71").Value = 3716.6365
$ws.Range("I71").Value = 1870.125
$ws.Range("J71").Value = 8640.666999999999
$ws.Range("K71").Value = 9350.625
$ws.Range("L71").Value = 43203.335
$ws.Range("M71").Value = -5606.625
$ws.Range("N71").Value = -50691.335

$ws.Range("H82").Value = 497.4375
$ws.Range("I82").Value = 523.8333
$ws.Range("J82").Value = 418.25
$ws.Range("K82").Value = 523.8333
$ws.Range("L82").Value = 418.25
$ws.Range("M82").Value = -162.8333
$ws.Range("N82").Value = -1140.25

$ws.Range("H85").Value = 497.4375
$ws.Range("I85").Value = 523.8333
$ws.Range("J85").Value = 418.25
$ws.Range("K85").Value = 523.8333
$ws.Range("L85").Value = 418.25
$ws.Range("M85").Value = 724.1667
$ws.Range("N85").Value = -2914.25

$ws.Range("H87").Value = 85555
$ws.Range("J87").Value = 85555
$ws.Range("L87").Value = 85555
$ws.Range("N87").Value = -87801

$ws.Range("H90").Value = 85555
$ws.Range("J90").Value = 85555
$ws.Range("L90").Value = 256665
$ws.Range("N90").Value = -267897

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = $null
$ws.Range("N111").Value = $null

$ws.Range("H122").Value = 4549.8335
$ws.Range("I122").Value = 4434.8
$ws.Range("J122").Value = 5125
$ws.Range("K122").Value = 13304.4
$ws.Range("L122").Value = 15375
$ws.Range("M122").Value = -10854.4
$ws.Range("N122").Value = -20275

$ws.Range("H126").Value = 5062.6665
$ws.Range("I126").Value = 4875.2
$ws.Range("K126").Value = 14625.6
$ws.Range("M126").Value = -12155.6

$ws.Range("H132").Value = 10006204
$ws.Range("I132").Value = 12506191
$ws.Range("J132").Value = 6258.8
$ws.Range("K132").Value = 37518573
$ws.Range("L132").Value = 18776.4
$ws.Range("M132").Value = -37516043
$ws.Range("N132").Value = -23836.4

$ws.Range("H136").Value = 2860.9048
$ws.Range("I136").Value = 3608.4167
$ws.Range("J136").Value = 1864.2222
$ws.Range("K136").Value = 10825.2501
$ws.Range("L136").Value = 5592.6666
$ws.Range("M136").Value = -8275.250100000001
$ws.Range("N136").Value = -10692.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 100001170
$ws.Range("J4").Value = 100001170
$ws.Range("L4").Value = 100001170
$ws.Range("N4").Value = -100001396

$ws.Range("H62").Value = 6003
$ws.Range("I62").Value = 4388.8
$ws.Range("J62").Value = 7011.875
$ws.Range("K62").Value = 4388.8
$ws.Range("L62").Value = 7011.875
$ws.Range("M62").Value = -3764.8
$ws.Range("N62").Value = -8259.875

$ws.Range("H65").Value = 6003
$ws.Range("I65").Value = 4388.8
$ws.Range("J65").Value = 7011.875
$ws.Range("K65").Value = 21944
$ws.Range("L65").Value = 35059.375
$ws.Range("M65").Value = -18824
$ws.Range("N65").Value = -41299.375

$ws.Range("H81").Value = 1821.4286
$ws.Range("I81").Value = 1952.4
$ws.Range("J81").Value = 1494
$ws.Range("K81").Value = 3904.8
$ws.Range("L81").Value = 2988
$ws.Range("M81").Value = -2843.8
$ws.Range("N81").Value = -5110

$ws.Range("H84").Value = 1821.4286
$ws.Range("I84").Value = 1952.4
$ws.Range("J84").Value = 1494
$ws.Range("K84").Value = 19524
$ws.Range("L84").Value = 14940
$ws.Range("M84").Value = -14220
$ws.Range("N84").Value = -25548

$ws.Range("H98").Value = 26496.334
$ws.Range("J98").Value = 26496.334
$ws.Range("L98").Value = 26496.334
$ws.Range("N98").Value = -32486.334

$ws.Range("H107").Value = 367.6875
$ws.Range("I107").Value = 419.5
$ws.Range("J107").Value = 212.25
$ws.Range("K107").Value = 1258.5
$ws.Range("L107").Value = 636.75
$ws.Range("M107").Value = 661.5
$ws.Range("N107").Value = -4476.75

$ws.Range("H122").Value = 1230.5834
$ws.Range("I122").Value = 1265.2727
$ws.Range("J122").Value = 849
$ws.Range("K122").Value = 3795.8181
$ws.Range("L122").Value = 2547
$ws.Range("M122").Value = -1345.8181
$ws.Range("N122").Value = -7447

$ws.Range("H132").Value = 500000000
$ws.Range("I132").Value = 500000000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1500000000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1499997470
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 41668130
$ws.Range("I136").Value = 45455956
$ws.Range("K136").Value = 136367868
$ws.Range("M136").Value = -136365318

